$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, pushing existing rows 3-8 down to 4-9
$ws.Rows.Item(3).Insert()

# Copy constant columns (A, B, C, E, F, G, H, I, J) from row 4 (the row that used to be row 3)
$ws.Range("A3").Value = $ws.Range("A4").Value2
$ws.Range("B3").Value = $ws.Range("B4").Value2
$ws.Range("C3").Value = $ws.Range("C4").Value2
$ws.Range("E3").Value = $ws.Range("E4").Value2
$ws.Range("F3").Value = $ws.Range("F4").Value2
$ws.Range("G3").Value = $ws.Range("G4").Value2
$ws.Range("H3").Value = $ws.Range("H4").Value2
$ws.Range("I3").Value = $ws.Range("I4").Value2
$ws.Range("J3").Value = $ws.Range("J4").Value2

# New row 3 data
$ws.Range("D3").Value = 44557
$ws.Range("D3").NumberFormat = $ws.Range("D4").NumberFormat
$ws.Range("K3").Value = "Lapins"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 250
$ws.Range("N3").Value = 9000
$ws.Range("O3").Value = 10000
$ws.Range("P3").Value = 9500
$ws.Range("Q3").Value = "`$/bandeja 10 kilos"
$ws.Range("R3").Value = "Provincia de Curicó"
$ws.Range("S3").Value = 950
$ws.Range("T3").Value = 10
